# GitHub Actions symbol-list refresh for the cryptos sheet.
# - Column D ("Price") holds numeric-looking text (e.g. "236.28"), not real
#   numbers, so each touched cell is switched to Text format ("@") before the
#   new value is assigned. That preserves exact formatting (trailing zeros,
#   no scientific notation for tiny values) instead of Excel auto-converting
#   the assignment into a Double.
# - A few rows were re-ranked (coins shifted to a different rank row, two
#   rows swapped order) which shows up as whole-row content changes across
#   B (Coin), C (Link), D (Price) and E (Volume(1h)) for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column updates (stored as text) ---------------------------------
$priceUpdates = [ordered]@{
    "D2"  = "236.28"
    "D4"  = "5.409"
    "D5"  = "0.05632"
    "D6"  = "3.377"
    "D7"  = "6.483"
    "D8"  = "1.067"
    "D9"  = "0.7822"
    "D11" = "0.07412"
    "D12" = "0.03176"
    "D13" = "0.02958"
    "D14" = "0.09265"
    "D15" = "0.001663"
    "D16" = "3.254"
    "D17" = "0.04735"
    "D18" = "0.0005786"
    "D19" = "0.006208"
    "D20" = "0.005120"
    "D21" = "0.001051"
    "D22" = "0.0001500"
    "D23" = "3.920"
    "D24" = "2.146"
    "D40" = "0.04056"
    "D41" = "0.007021"
    "D42" = "0.1043"
    "D43" = "0.002631"
    "D44" = "0.009401"
    "D46" = "0.00000000750"
    "D47" = "0.6748"
    "D48" = "0.03992"
    "D49" = "0.00002099"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
}

# --- Coin / Link / Volume(1h) text updates (rank reshuffle + renames) ------
$textUpdates = [ordered]@{
    "B18" = "One"
    "C18" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E18" = "17OneONE"

    "B19" = "TigerCash"
    "C19" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "E19" = "18TigerCashTCH"

    "B20" = "HotbitToken"
    "C20" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
    "E20" = "19HotbitTokenHTB"

    "B21" = "BitKan"
    "C21" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
    "E21" = "20BitKanKAN"

    "B22" = "NitroEx"
    "C22" = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
    "E22" = "21NitroExNTX"

    "B23" = "LEO"
    "C23" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "E23" = "22LEOLEO"

    "B24" = "BTSEToken"
    "C24" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "E24" = "23BTSETokenBTSE"

    "E27" = "26UpBotsUBXTBestin24h"

    "B42" = "BKEXToken"
    "C42" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "E42" = "41BKEXTokenBKK"

    "B43" = "CEJI"
    "C43" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "E43" = "42CEJICEJI"
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = $textUpdates[$cellRef]
}
